$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1167.7593
$ws.Range("I137").Value = 889.6585
$ws.Range("J137").Value = 2044.8462
$ws.Range("K137").Value = 2668.9755
$ws.Range("L137").Value = 6134.5386
$ws.Range("M137").Value = -118.9755
$ws.Range("N137").Value = -11234.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1695.4
$ws.Range("I61").Value = 1119.75
$ws.Range("K61").Value = 1119.75
$ws.Range("M61").Value = -907.75
$ws.Range("H74").Value = 532.5714
$ws.Range("I74").Value = 538
$ws.Range("J74").Value = 234
$ws.Range("K74").Value = 538
$ws.Range("L74").Value = 234
$ws.Range("M74").Value = 336
$ws.Range("N74").Value = -1982
$ws.Range("H77").Value = 532.5714
$ws.Range("I77").Value = 538
$ws.Range("J77").Value = 234
$ws.Range("K77").Value = 2690
$ws.Range("L77").Value = 1170
$ws.Range("M77").Value = 1678
$ws.Range("N77").Value = -9906
$ws.Range("H136").Value = 1695.4
$ws.Range("I136").Value = 1119.75
$ws.Range("K136").Value = 3359.25
$ws.Range("M136").Value = -809.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 26501.244
$ws.Range("I134").Value = 35128.3
$ws.Range("J134").Value = 2972.9092
$ws.Range("K134").Value = 105384.9
$ws.Range("L134").Value = 8918.7276
$ws.Range("M134").Value = -102849.9
$ws.Range("N134").Value = -13988.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3032506
$ws.Range("I31").Value = 1534.683
$ws.Range("J31").Value = 11908922
$ws.Range("K31").Value = 1534.683
$ws.Range("L31").Value = 11908922
$ws.Range("M31").Value = -1239.683
$ws.Range("N31").Value = -11909512
$ws.Range("H34").Value = 3032506
$ws.Range("I34").Value = 1534.683
$ws.Range("J34").Value = 11908922
$ws.Range("K34").Value = 1534.683
$ws.Range("L34").Value = 11908922
$ws.Range("M34").Value = -1332.683
$ws.Range("N34").Value = -11909326
$ws.Range("H58").Value = 7408312.5
$ws.Range("I58").Value = 893.7879
$ws.Range("J58").Value = 27778714
$ws.Range("K58").Value = 893.7879
$ws.Range("L58").Value = 27778714
$ws.Range("M58").Value = -690.7879
$ws.Range("N58").Value = -27779120
$ws.Range("H132").Value = 1709.2333
$ws.Range("I132").Value = 1523.7142
$ws.Range("J132").Value = 2535.6365
$ws.Range("K132").Value = 4571.142599999999
$ws.Range("L132").Value = 7606.9095
$ws.Range("M132").Value = -2041.142599999999
$ws.Range("N132").Value = -12666.9095
$ws.Range("H134").Value = 1315.919
$ws.Range("I134").Value = 1271.5312
$ws.Range("J134").Value = 1600
$ws.Range("K134").Value = 3814.5936
$ws.Range("L134").Value = 4800
$ws.Range("M134").Value = -1279.5936
$ws.Range("N134").Value = -9870
$ws.Range("H136").Value = 7408312.5
$ws.Range("I136").Value = 893.7879
$ws.Range("J136").Value = 27778714
$ws.Range("K136").Value = 2681.3637
$ws.Range("L136").Value = 83336142
$ws.Range("M136").Value = -131.3636999999999
$ws.Range("N136").Value = -83341242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 2231.7058
$ws.Range("I123").Value = 694
$ws.Range("J123").Value = 4428.4287
$ws.Range("K123").Value = 2082
$ws.Range("L123").Value = 13285.2861
$ws.Range("M123").Value = 368
$ws.Range("N123").Value = -18185.2861
$ws.Range("H129").Value = 1250.7142
$ws.Range("I129").Value = 935.6667
$ws.Range("J129").Value = 1487
$ws.Range("K129").Value = 2807.0001
$ws.Range("L129").Value = 4461
$ws.Range("M129").Value = 2192.9999
$ws.Range("N129").Value = -14461
$ws.Range("H130").Value = 2387.2727
$ws.Range("I130").Value = 1208.5714
$ws.Range("J130").Value = 4450
$ws.Range("K130").Value = 3625.7142
$ws.Range("L130").Value = 13350
$ws.Range("M130").Value = 1394.2858
$ws.Range("N130").Value = -23390
$ws.Range("H131").Value = 4663504
$ws.Range("I131").Value = 18601.727
$ws.Range("J131").Value = 7352657.5
$ws.Range("K131").Value = 55805.181
$ws.Range("L131").Value = 22057972.5
$ws.Range("M131").Value = -50765.181
$ws.Range("N131").Value = -22068052.5
$ws.Range("H133").Value = 2450.4546
$ws.Range("I133").Value = 2435.8823
$ws.Range("K133").Value = 7307.646900000001
$ws.Range("M133").Value = -2247.646900000001
$ws.Range("H134").Value = 2450.7693
$ws.Range("I134").Value = 2420
$ws.Range("K134").Value = 7260
$ws.Range("M134").Value = -2190
$ws.Range("H136").Value = 1874.7142
$ws.Range("I136").Value = 1197.7778
$ws.Range("J136").Value = 3093.2
$ws.Range("K136").Value = 3593.3334
$ws.Range("L136").Value = 9279.599999999999
$ws.Range("M136").Value = 1506.6666
$ws.Range("N136").Value = -19479.6
$ws.Range("H137").Value = 63197740
$ws.Range("I137").Value = 41681796
$ws.Range("J137").Value = 75492570
$ws.Range("K137").Value = 125045388
$ws.Range("L137").Value = 226477710
$ws.Range("M137").Value = -125040288
$ws.Range("N137").Value = -226487910
$ws.Range("H138").Value = 2082.8572
$ws.Range("I138").Value = 1896.6666
$ws.Range("K138").Value = 5689.9998
$ws.Range("M138").Value = -549.9997999999996
$ws.Range("H139").Value = 2296
$ws.Range("I139").Value = 1670
$ws.Range("K139").Value = 5010
$ws.Range("M139").Value = 130
$ws.Range("H140").Value = 2016.9231
$ws.Range("I140").Value = 1810
$ws.Range("K140").Value = 5430
$ws.Range("M140").Value = -250
$ws.Range("H141").Value = 2377.4333
$ws.Range("I141").Value = 2281.6538
$ws.Range("K141").Value = 6844.9614
$ws.Range("M141").Value = -1664.9614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12347.286
$ws.Range("I132").Value = 16864.715
$ws.Range("J132").Value = 3312.4285
$ws.Range("K132").Value = 50594.145
$ws.Range("L132").Value = 9937.2855
$ws.Range("M132").Value = -48064.145
$ws.Range("N132").Value = -14997.2855
$ws.Range("H136").Value = 2790.709
$ws.Range("I136").Value = 2760.1875
$ws.Range("K136").Value = 8280.5625
$ws.Range("M136").Value = -5730.5625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1516.4468
$ws.Range("I132").Value = 1241.9
$ws.Range("J132").Value = 3085.2856
$ws.Range("K132").Value = 3725.7
$ws.Range("L132").Value = 9255.856800000001
$ws.Range("M132").Value = -1195.7
$ws.Range("N132").Value = -14315.8568
$ws.Range("H136").Value = 2964.1187
$ws.Range("I136").Value = 3307.1914
$ws.Range("J136").Value = 1620.4166
$ws.Range("K136").Value = 9921.574200000001
$ws.Range("L136").Value = 4861.2498
$ws.Range("M136").Value = -7371.574200000001
$ws.Range("N136").Value = -9961.2498
